# Updates the "Volume(1h)" percentage column (column E) on the active
# worksheet to the refreshed crypto price-change values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -2.74%  "
$ws.Range("E3").Value = "  -2.98%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("E5").Value = "  -3.12%  "
$ws.Range("E6").Value = "  -4.13%  "
$ws.Range("E7").Value = "  +0.08%  "
$ws.Range("E8").Value = "  -2.93%  "
$ws.Range("E9").Value = "  -2.30%  "
$ws.Range("E10").Value = "  -1.15%  "
$ws.Range("E11").Value = "  -3.45%  "
$ws.Range("E12").Value = "  -1.87%  "
$ws.Range("E13").Value = "  -2.91%  "
$ws.Range("E14").Value = "  +0.26%  "
$ws.Range("E15").Value = "  -5.67%  "
$ws.Range("E16").Value = "  -2.98%  "
$ws.Range("E17").Value = "  -3.21%  "
$ws.Range("E18").Value = "  -2.75%  "
$ws.Range("E19").Value = "  -3.62%  "
$ws.Range("E20").Value = "  -0.93%  "
$ws.Range("E21").Value = "  -3.47%  "
$ws.Range("E22").Value = "  -1.41%  "
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("E24").Value = "  -4.24%  "
$ws.Range("E25").Value = "  -0.04%  "
$ws.Range("E26").Value = "  -3.49%  "
$ws.Range("E27").Value = "  -7.61%  "
$ws.Range("E28").Value = "  -4.22%  "
$ws.Range("E29").Value = "  +0.06%  "
$ws.Range("E30").Value = "  -5.44%  "
$ws.Range("E31").Value = "  -0.01%  "
$ws.Range("E33").Value = "  -4.25%  "
$ws.Range("E34").Value = "  -2.09%  "
$ws.Range("E35").Value = "  -5.41%  "
$ws.Range("E36").Value = "  -5.69%  "
$ws.Range("E37").Value = "  -2.16%  "
$ws.Range("E38").Value = "  -5.15%  "
$ws.Range("E39").Value = "  -2.39%  "
$ws.Range("E40").Value = "  -3.07%  "
$ws.Range("E41").Value = "  -11.83%  "
$ws.Range("E42").Value = "  -5.46%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("E44").Value = "  -3.62%  "
$ws.Range("E45").Value = "  -3.87%  "
$ws.Range("E46").Value = "  -4.84%  "
$ws.Range("E47").Value = "  -4.22%  "
$ws.Range("E48").Value = "  -7.10%  "
$ws.Range("E49").Value = "  -0.21%  "
$ws.Range("E50").Value = "  -5.45%  "
$ws.Range("E51").Value = "  -6.38%  "
